# Roll the yearly columns forward by one fiscal year:
#  - drop the oldest year ("...1396/12") from the header row
#  - shift each data row's values left by one column (E<-F, F<-G, G<-H, H<-I)
#  - populate the newly-freed last column (I) with the new fiscal year's figures
#  - add the new header label ("...1401/12") in column I

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows (year labels) ---------------------------------------
# Row 8 and Row 24 both show the same 5 year headers across E:I.
$newHeader = "دوازده ماهه منتهی به 1401/12"

foreach ($headerRow in 8, 24) {
    $shifted = $ws.Range("F" + $headerRow + ":I" + $headerRow).Value()
    $ws.Range("E" + $headerRow + ":H" + $headerRow).Value = $shifted
    $ws.Range("I" + $headerRow).Value = $newHeader
}

# --- Data rows ---------------------------------------------------------
# For each data row, shift E:I left by one column (drop old E, keep F..I
# in E..H), then put the new fiscal year's value into I.

$newLastValues = @{
    13 = 16619
    14 = 36627
    15 = 8463
    16 = 11398
    17 = 506587
    19 = 129109
    20 = 708803
    26 = 400
    27 = 210
}

foreach ($row in $newLastValues.Keys) {
    $shifted = $ws.Range("F" + $row + ":I" + $row).Value()
    $ws.Range("E" + $row + ":H" + $row).Value = $shifted
    $ws.Range("I" + $row).Value = $newLastValues[$row]
}
